# Rename sheet "DimensionxFactorxPregunta" -> "Dimensions",
# make it the active sheet/tab (was "Hypotheses"), and update its
# selected cell from H13 to G17.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DimensionxFactorxPregunta")
$ws.Name = "Dimensions"

# Make the renamed sheet the active one (updates workbookView activeTab
# and moves tabSelected="1" from the previously active sheet to this one).
$ws.Activate()

# Update the selection/active cell on this sheet.
$ws.Range("G17").Select()
